$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(349.46829586836225, 0.9703288855936801, -0.9797427175039001),
    @(355.26867350781674, 1.0356920359203379, -1.0422382482031072),
    @(364.81360632880165, 0.886195587501846, -1.102525546796564),
    @(336.0923025285776, 1.1219619154265323, -0.9023655792136507),
    @(353.0957623931573, 1.0473903084237475, -1.0921507382820521),
    @(347.70425452094946, 0.779046011731997, -0.804996199739501),
    @(379.8052079409878, 0.7016762681413808, -1.264789233063625),
    @(361.44259169021694, 0.8783641774035735, -1.0816614299809972),
    @(346.85630286365466, 0.9256261399502125, -0.9871272144856588),
    @(339.39908595113366, 0.9641736251853227, -0.861056924050909),
    @(349.4611760784422, 1.004488626224144, -0.9176523685807214),
    @(332.5049845275216, 1.113500580078438, -0.7598845663538311),
    @(369.1932244523061, 0.6663679371714449, -1.0545946586781179),
    @(339.70380017894377, 1.166851312434524, -0.9381441103267879),
    @(349.90045636202456, 0.9720347123644691, -1.0429422493660254),
    @(340.62604571770515, 0.9834710580960611, -0.8010119955256351),
    @(322.0170346054832, 1.4557202053343814, -0.9425764728223761),
    @(340.08688845224555, 1.072284842231723, -0.8633942508650241),
    @(327.2282511162636, 1.2393283915745281, -0.6919790460026949),
    @(364.3536657162726, 0.8423786052169173, -1.0983339920409605),
    @(342.468042196242, 1.1793065456100336, -1.1190438932089033),
    @(335.0967149658947, 0.8593091793628427, -0.8105394126987004),
    @(355.14699999311495, 0.9239026558554356, -1.1503061603341214),
    @(345.0075832757741, 0.9607028410605132, -0.871569514534028),
    @(340.6144351923502, 1.1198636616828979, -0.8462195974498432),
    @(335.9001511062851, 1.080756309162901, -0.7161581308479826),
    @(351.7982056002291, 1.0202022577106502, -0.9240456373012649),
    @(349.4624545428592, 0.9700329616851953, -1.1403042210223542),
    @(360.3587072918426, 0.8567552521019224, -1.0846685474113151),
    @(348.01947618698784, 0.768577694891657, -0.992605795877557),
    @(347.862663802525, 0.9412665671584767, -0.9038411727779017),
    @(345.3130812487108, 0.8445193534464887, -0.6585094318424888),
    @(349.00038800492047, 0.9448862820408439, -0.8443265119272441),
    @(339.61608580025705, 1.2441677731464307, -0.9261257181301047),
    @(374.818348956407, 0.8895796635287818, -1.1755614171275792),
    @(362.7421212079232, 0.8755170721429469, -0.946085867889396),
    @(344.18662295529737, 1.2339789766123992, -1.0683721601456604),
    @(355.6445571975587, 0.7548667501783878, -0.9423755908733323),
    @(332.0142642878737, 1.1487933238379722, -0.8599360297777824),
    @(372.07382948042493, 0.7126316801831322, -1.2143689857850954),
    @(358.9201124665608, 1.0500057965822802, -1.2495088581836515),
    @(350.8887166063648, 0.7892375512286196, -0.94984517447213),
    @(340.5310733598086, 0.9990251176926074, -0.8259807334069081),
    @(345.2863191470257, 0.9684467866529315, -0.8864026387546545),
    @(337.0583044987576, 0.9516234797608378, -0.8565075668608849),
    @(359.6306648658198, 0.9715653959283321, -1.1574398306303986),
    @(367.6249325946656, 0.8207103540867329, -1.0382021911367594),
    @(350.5096820946877, 1.1391644845679647, -1.1236465700333156),
    @(347.13061576404596, 1.1879318214853998, -1.0625580211109227),
    @(339.0306565113941, 0.998751633674858, -0.7738368464669888)
)

$values = New-Object 'object[,]' 50,3
for ($i = 0; $i -lt 50; $i++) {
    for ($j = 0; $j -lt 3; $j++) {
        $values[$i, $j] = $data[$i][$j]
    }
}

$ws.Range("A2:C51").Value = $values
